$d = $word.ActiveDocument

# 1) "Fecha de la solicitud" -> date cell: 19/02/2025 -> 24/02/2025
#    (Table 1, row 6, col 2)
$t1 = $d.Tables.Item(1)
$t1.Rows.Item(6).Cells.Item(2).Range.Text = "24/02/2025"

# 2) "Se encuentra matriculado o en reserva de cupo..." NO -> SI
#    (Table 2, row 3, col 2)
$t2 = $d.Tables.Item(2)
$t2.Rows.Item(3).Cells.Item(2).Range.Text = "SI"

# 3) "Componente de LIBRE ELECCIÓN" header row: second Código/Asignatura pair
#    gets the _CC suffix (Table 6, row 1, cols 4 and 5)
$t6 = $d.Tables.Item(6)
$t6.Rows.Item(1).Cells.Item(4).Range.Text = "Código_CC"
$t6.Rows.Item(1).Cells.Item(5).Range.Text = "Asignatura_CC"

# 4) Section title text change
$d.Content.Find.Execute("Pendientes en Fundamentación - Optativa T", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Pendientes en Disciplinar - Optativa T", 2)
